# Add a "Save" column (H) to the s_vals sheet, matching the style of the
# other header cells, with 0/1 flag values for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: same look as the other header cells (B1:G1) -> copy format
# from G1 ("sum") onto H1, then set its text.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for the new column.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
